$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 658.087
$ws.Range("I6").Value = 532.35
$ws.Range("K6").Value = 1597.05
$ws.Range("M6").Value = -1485.05
$ws.Range("H12").Value = 788.3333
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 824.375
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 824.375
$ws.Range("M12").Value = -330
$ws.Range("N12").Value = -1164.375
$ws.Range("H40").Value = 7713.8096
$ws.Range("I40").Value = 12657.777
$ws.Range("J40").Value = 4005.8333
$ws.Range("K40").Value = 12657.777
$ws.Range("L40").Value = 4005.8333
$ws.Range("M40").Value = -12482.777
$ws.Range("N40").Value = -4355.8333
$ws.Range("H43").Value = 1536.5714
$ws.Range("I43").Value = 1200
$ws.Range("J43").Value = 1628.3636
$ws.Range("K43").Value = 1200
$ws.Range("L43").Value = 1628.3636
$ws.Range("M43").Value = -1131
$ws.Range("N43").Value = -1766.3636
$ws.Range("H111").Value = 13895046
$ws.Range("I111").Value = 22230974
$ws.Range("K111").Value = 66692922
$ws.Range("M111").Value = -66689855
$ws.Range("H113").Value = 8070.909
$ws.Range("I113").Value = 7666.6665
$ws.Range("K113").Value = 7666.6665
$ws.Range("M113").Value = -4412.6665
$ws.Range("H116").Value = 6011
$ws.Range("J116").Value = 6011
$ws.Range("L116").Value = 6011
$ws.Range("N116").Value = -12895
$ws.Range("H121").Value = 2505.476
$ws.Range("J121").Value = 2530.75
$ws.Range("L121").Value = 7592.25
$ws.Range("N121").Value = -11086.25
$ws.Range("H127").Value = 927.63635
$ws.Range("I127").Value = 927.63635
$ws.Range("K127").Value = 2782.90905
$ws.Range("M127").Value = 2177.09095
$ws.Range("H137").Value = 78520.25999999999
$ws.Range("I137").Value = 178746.9
$ws.Range("J137").Value = 1422.8462
$ws.Range("K137").Value = 536240.7
$ws.Range("L137").Value = 4268.5386
$ws.Range("M137").Value = -533690.7
$ws.Range("N137").Value = -9368.5386
$ws.Range("H138").Value = 4783.5835
$ws.Range("I138").Value = 3968.1
$ws.Range("J138").Value = 4946.68
$ws.Range("K138").Value = 11904.3
$ws.Range("L138").Value = 14840.04
$ws.Range("M138").Value = -6764.299999999999
$ws.Range("N138").Value = -25120.04
$ws.Range("H141").Value = 6060.0454
$ws.Range("I141").Value = 6743.2104
$ws.Range("J141").Value = 1733.3334
$ws.Range("K141").Value = 20229.6312
$ws.Range("L141").Value = 5200.0002
$ws.Range("M141").Value = -15049.6312
$ws.Range("N141").Value = -15560.0002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11113509
$ws.Range("I2").Value = 18520184
$ws.Range("J2").Value = 3497.5
$ws.Range("K2").Value = 18520184
$ws.Range("L2").Value = 3497.5
$ws.Range("M2").Value = -18520071
$ws.Range("N2").Value = -3723.5
$ws.Range("H32").Value = 2768.244
$ws.Range("I32").Value = 1693.1096
$ws.Range("K32").Value = 1693.1096
$ws.Range("M32").Value = -1406.1096
$ws.Range("H55").Value = 31015
$ws.Range("I55").Value = 3048
$ws.Range("J55").Value = 44998.5
$ws.Range("K55").Value = 3048
$ws.Range("L55").Value = 44998.5
$ws.Range("M55").Value = -2733
$ws.Range("N55").Value = -45628.5
$ws.Range("H60").Value = 10525.25
$ws.Range("I60").Value = 6050.5
$ws.Range("J60").Value = 15000
$ws.Range("K60").Value = 6050.5
$ws.Range("L60").Value = 15000
$ws.Range("M60").Value = -5317.5
$ws.Range("N60").Value = -16466
$ws.Range("H61").Value = 8903.611000000001
$ws.Range("I61").Value = 9133.529
$ws.Range("K61").Value = 9133.529
$ws.Range("M61").Value = -8921.529
$ws.Range("H74").Value = 51293.645
$ws.Range("I74").Value = 9891.281000000001
$ws.Range("J74").Value = 183781.2
$ws.Range("K74").Value = 9891.281000000001
$ws.Range("L74").Value = 183781.2
$ws.Range("M74").Value = -9017.281000000001
$ws.Range("N74").Value = -185529.2
$ws.Range("H77").Value = 51293.645
$ws.Range("I77").Value = 9891.281000000001
$ws.Range("J77").Value = 183781.2
$ws.Range("K77").Value = 49456.40500000001
$ws.Range("L77").Value = 918906
$ws.Range("M77").Value = -45088.40500000001
$ws.Range("N77").Value = -927642
$ws.Range("H116").Value = 11113509
$ws.Range("I116").Value = 18520184
$ws.Range("J116").Value = 3497.5
$ws.Range("K116").Value = 18520184
$ws.Range("L116").Value = 3497.5
$ws.Range("M116").Value = -18517890
$ws.Range("N116").Value = -8085.5
$ws.Range("H122").Value = 427614.03
$ws.Range("I122").Value = 2128.9714
$ws.Range("J122").Value = 1491326.6
$ws.Range("K122").Value = 6386.914199999999
$ws.Range("L122").Value = 4473979.800000001
$ws.Range("M122").Value = -3936.914199999999
$ws.Range("N122").Value = -4478879.800000001
$ws.Range("H132").Value = 5346.1626
$ws.Range("J132").Value = 4882.4116
$ws.Range("L132").Value = 14647.2348
$ws.Range("N132").Value = -19707.2348
$ws.Range("H136").Value = 8903.611000000001
$ws.Range("I136").Value = 9133.529
$ws.Range("K136").Value = 27400.587
$ws.Range("M136").Value = -24850.587

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11113509
$ws.Range("I3").Value = 18520184
$ws.Range("J3").Value = 3497.5
$ws.Range("K3").Value = 18520184
$ws.Range("L3").Value = 3497.5
$ws.Range("M3").Value = -18520070
$ws.Range("N3").Value = -3725.5
$ws.Range("H29").Value = 154249.75
$ws.Range("I29").Value = 154249.75
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 154249.75
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -153960.75
$ws.Range("N29").ClearContents()
$ws.Range("H36").Value = 1931.3334
$ws.Range("I36").Value = 1931.3334
$ws.Range("K36").Value = 1931.3334
$ws.Range("M36").Value = -1397.3334
$ws.Range("H86").Value = 6265346
$ws.Range("I86").Value = 9111767
$ws.Range("K86").Value = 9111767
$ws.Range("M86").Value = -9110644
$ws.Range("H89").Value = 6265346
$ws.Range("I89").Value = 9111767
$ws.Range("K89").Value = 45558835
$ws.Range("M89").Value = -45553219
$ws.Range("H99").Value = 6214279
$ws.Range("I99").Value = 9526635
$ws.Range("K99").Value = 9526635
$ws.Range("M99").Value = -9525137
$ws.Range("H105").Value = 2887138.5
$ws.Range("I105").Value = 3175752.2
$ws.Range("K105").Value = 3175752.2
$ws.Range("M105").Value = -3174005.2
$ws.Range("H134").Value = 5516.0513
$ws.Range("I134").Value = 4397.9697
$ws.Range("K134").Value = 13193.9091
$ws.Range("M134").Value = -10658.9091

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1087.9412
$ws.Range("I4").Value = 991.3333
$ws.Range("K4").Value = 991.3333
$ws.Range("M4").Value = -879.3333
$ws.Range("H7").Value = 230.34782
$ws.Range("I7").Value = 84.46666999999999
$ws.Range("K7").Value = 84.46666999999999
$ws.Range("M7").Value = 28.53333000000001
$ws.Range("H31").Value = 4091.97
$ws.Range("I31").Value = 1520
$ws.Range("J31").Value = 4618.7593
$ws.Range("K31").Value = 1520
$ws.Range("L31").Value = 4618.7593
$ws.Range("M31").Value = -1225
$ws.Range("N31").Value = -5208.7593
$ws.Range("H34").Value = 4091.97
$ws.Range("I34").Value = 1520
$ws.Range("J34").Value = 4618.7593
$ws.Range("K34").Value = 1520
$ws.Range("L34").Value = 4618.7593
$ws.Range("M34").Value = -1318
$ws.Range("N34").Value = -5022.7593
$ws.Range("H58").Value = 2845.7896
$ws.Range("I58").Value = 2677.6155
$ws.Range("J58").Value = 3210.1667
$ws.Range("K58").Value = 2677.6155
$ws.Range("L58").Value = 3210.1667
$ws.Range("M58").Value = -2474.6155
$ws.Range("N58").Value = -3616.1667
$ws.Range("H59").Value = 40299.8
$ws.Range("J59").Value = 47499.668
$ws.Range("L59").Value = 47499.668
$ws.Range("N59").Value = -49789.668
$ws.Range("H94").Value = 5980.4707
$ws.Range("J94").Value = 6182.8335
$ws.Range("L94").Value = 6182.8335
$ws.Range("N94").Value = -7084.8335
$ws.Range("H99").Value = 3147.625
$ws.Range("I99").Value = 2892.2856
$ws.Range("K99").Value = 2892.2856
$ws.Range("M99").Value = -1394.2856
$ws.Range("H126").Value = 3147.625
$ws.Range("I126").Value = 2892.2856
$ws.Range("K126").Value = 8676.856800000001
$ws.Range("M126").Value = -6206.856800000001
$ws.Range("H132").Value = 93593.73
$ws.Range("I132").Value = 145315.14
$ws.Range("K132").Value = 435945.42
$ws.Range("M132").Value = -433415.42
$ws.Range("H136").Value = 2845.7896
$ws.Range("I136").Value = 2677.6155
$ws.Range("J136").Value = 3210.1667
$ws.Range("K136").Value = 8032.8465
$ws.Range("L136").Value = 9630.500100000001
$ws.Range("M136").Value = -5482.8465
$ws.Range("N136").Value = -14730.5001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 173.03922
$ws.Range("I2").Value = 96.36
$ws.Range("J2").Value = 246.76923
$ws.Range("K2").Value = 578.16
$ws.Range("L2").Value = 1480.61538
$ws.Range("M2").Value = -465.16
$ws.Range("N2").Value = -1706.61538
$ws.Range("H5").Value = 4821.913
$ws.Range("I5").Value = 439.75
$ws.Range("K5").Value = 1319.25
$ws.Range("M5").Value = -1207.25
$ws.Range("H6").Value = 407.625
$ws.Range("I6").Value = 414.66666
$ws.Range("J6").Value = 302
$ws.Range("K6").Value = 1243.99998
$ws.Range("L6").Value = 906
$ws.Range("M6").Value = -1130.99998
$ws.Range("N6").Value = -1132
$ws.Range("H7").Value = 1261
$ws.Range("I7").Value = 1825.5834
$ws.Range("K7").Value = 5476.7502
$ws.Range("M7").Value = -5364.7502
$ws.Range("H11").Value = 2267
$ws.Range("I11").Value = 1334.2858
$ws.Range("J11").Value = 3454.0908
$ws.Range("K11").Value = 4002.8574
$ws.Range("L11").Value = 10362.2724
$ws.Range("M11").Value = -3862.8574
$ws.Range("N11").Value = -10642.2724
$ws.Range("H12").Value = 48169.316
$ws.Range("I12").Value = 111121.625
$ws.Range("K12").Value = 333364.875
$ws.Range("M12").Value = -333191.875
$ws.Range("H29").Value = 111.8
$ws.Range("I29").Value = 139.66667
$ws.Range("J29").Value = 99.85714
$ws.Range("K29").Value = 419.00001
$ws.Range("L29").Value = 299.57142
$ws.Range("M29").Value = -142.00001
$ws.Range("N29").Value = -853.57142
$ws.Range("H40").Value = 21.947369
$ws.Range("J40").Value = 24.8
$ws.Range("L40").Value = 99.2
$ws.Range("N40").Value = -237.2
$ws.Range("H68").Value = 1660.8572
$ws.Range("I68").Value = 1445.1052
$ws.Range("J68").Value = 2116.3333
$ws.Range("K68").Value = 4335.3156
$ws.Range("L68").Value = 6348.999899999999
$ws.Range("M68").Value = -3524.3156
$ws.Range("N68").Value = -7970.999899999999
$ws.Range("H71").Value = 1660.8572
$ws.Range("I71").Value = 1445.1052
$ws.Range("J71").Value = 2116.3333
$ws.Range("K71").Value = 13005.9468
$ws.Range("L71").Value = 19046.9997
$ws.Range("M71").Value = -8949.9468
$ws.Range("N71").Value = -27158.9997
$ws.Range("H104").Value = 2030.8
$ws.Range("I104").Value = 226
$ws.Range("J104").Value = 2482
$ws.Range("K104").Value = 678
$ws.Range("L104").Value = 7446
$ws.Range("M104").Value = 1943
$ws.Range("N104").Value = -12688
$ws.Range("H135").Value = 4821.913
$ws.Range("I135").Value = 439.75
$ws.Range("K135").Value = 3957.75
$ws.Range("M135").Value = -1422.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3000
$ws.Range("I14").Value = 1000
$ws.Range("K14").Value = 1000
$ws.Range("M14").Value = -832
$ws.Range("H94").Value = 22500
$ws.Range("J94").Value = 22500
$ws.Range("L94").Value = 22500
$ws.Range("N94").Value = -23852
$ws.Range("H126").Value = 8053110
$ws.Range("I126").Value = 7578722.5
$ws.Range("K126").Value = 22736167.5
$ws.Range("M126").Value = -22733697.5
$ws.Range("H132").Value = 6782.9062
$ws.Range("I132").Value = 3787.2964
$ws.Range("J132").Value = 22959.2
$ws.Range("K132").Value = 11361.8892
$ws.Range("L132").Value = 68877.60000000001
$ws.Range("M132").Value = -8831.889200000001
$ws.Range("N132").Value = -73937.60000000001
$ws.Range("H134").Value = 33883.668
$ws.Range("J134").Value = 33883.668
$ws.Range("L134").Value = 101651.004
$ws.Range("N134").Value = -106721.004
$ws.Range("H135").Value = 70555.55499999999
$ws.Range("J135").Value = 70555.55499999999
$ws.Range("L135").Value = 70555.55499999999
$ws.Range("N135").Value = -80695.55499999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10000
$ws.Range("I2").Value = 5000
$ws.Range("K2").Value = 5000
$ws.Range("M2").Value = -4888
$ws.Range("H68").Value = 3805.125
$ws.Range("I68").Value = 3063
$ws.Range("K68").Value = 3063
$ws.Range("M68").Value = -2314
$ws.Range("H71").Value = 3805.125
$ws.Range("I71").Value = 3063
$ws.Range("K71").Value = 15315
$ws.Range("M71").Value = -11571
$ws.Range("H93").Value = 15876180
$ws.Range("I93").Value = 25642750
$ws.Range("J93").Value = 5504
$ws.Range("K93").Value = 25642750
$ws.Range("L93").Value = 5504
$ws.Range("M93").Value = -25641502
$ws.Range("N93").Value = -8000
$ws.Range("H132").Value = 4054.3333
$ws.Range("I132").Value = 3805.7693
$ws.Range("J132").Value = 4700.6
$ws.Range("K132").Value = 11417.3079
$ws.Range("L132").Value = 14101.8
$ws.Range("M132").Value = -8887.3079
$ws.Range("N132").Value = -19161.8
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10011.143
$ws.Range("J45").Value = 10011.143
$ws.Range("L45").Value = 10011.143
$ws.Range("N45").Value = -10993.143
$ws.Range("H126").Value = 3571.625
$ws.Range("I126").Value = 3207.0356
$ws.Range("K126").Value = 9621.106800000001
$ws.Range("M126").Value = -7151.106800000001
$ws.Range("H132").Value = 19827202
$ws.Range("I132").Value = 23813962
$ws.Range("J132").Value = 1222322.1
$ws.Range("K132").Value = 71441886
$ws.Range("L132").Value = 3666966.3
$ws.Range("M132").Value = -71439356
$ws.Range("N132").Value = -3672026.3
$ws.Range("H136").Value = 5048.15
$ws.Range("I136").Value = 6169.476
$ws.Range("J136").Value = 2431.7222
$ws.Range("K136").Value = 18508.428
$ws.Range("L136").Value = 7295.1666
$ws.Range("M136").Value = -15958.428
$ws.Range("N136").Value = -12395.1666
$ws.Range("H139").Value = 30476.666
$ws.Range("J139").Value = 30476.666
$ws.Range("L139").Value = 30476.666
$ws.Range("N139").Value = -40756.666
$ws.Range("H140").Value = 22498.25
$ws.Range("J140").Value = 22498.25
$ws.Range("L140").Value = 22498.25
$ws.Range("N140").Value = -32858.25
